# Update loading_percent values for Case_1_10 (380 kV case)
# Commit: "case with 380 kV done"
# Writes new computed loading-percent results into the existing results grid
# (rows 2-25, columns B/D/E/F/G/H/I/J/L/N). Columns C/K/M remain 0 (unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 17.97725937601226
$ws.Range("D2").Value = 8.149424081277266
$ws.Range("E2").Value = 13.04721537596023
$ws.Range("F2").Value = 37.30377486676422
$ws.Range("G2").Value = 45.37427097910122
$ws.Range("H2").Value = 17.79196693951939
$ws.Range("I2").Value = 25.38690279205079
$ws.Range("J2").Value = 9.91454894342046
$ws.Range("L2").Value = 13.43837388620201
$ws.Range("N2").Value = 18.37422271430653

# Row 3
$ws.Range("B3").Value = 17.65268489091354
$ws.Range("D3").Value = 8.066555722567731
$ws.Range("E3").Value = 12.9003322832672
$ws.Range("F3").Value = 37.20481862321027
$ws.Range("G3").Value = 44.8485396325957
$ws.Range("H3").Value = 17.77857222290829
$ws.Range("I3").Value = 25.54584260771767
$ws.Range("J3").Value = 9.867384831987568
$ws.Range("L3").Value = 13.2401757425236
$ws.Range("N3").Value = 18.41548558739757

# Row 4
$ws.Range("B4").Value = 17.45315289488404
$ws.Range("D4").Value = 8.01447824721417
$ws.Range("E4").Value = 12.80834658710939
$ws.Range("F4").Value = 37.15577890128522
$ws.Range("G4").Value = 44.5403741767387
$ws.Range("H4").Value = 17.77509798603262
$ws.Range("I4").Value = 25.64868481696153
$ws.Range("J4").Value = 9.838102571993637
$ws.Range("L4").Value = 13.11939593948662
$ws.Range("N4").Value = 18.44288460895341

# Row 5
$ws.Range("B5").Value = 17.37188473315764
$ws.Range("D5").Value = 7.992962597681598
$ws.Range("E5").Value = 12.77042751142811
$ws.Range("F5").Value = 37.13874962848143
$ws.Range("G5").Value = 44.41862309677484
$ws.Range("H5").Value = 17.77487592687035
$ws.Range("I5").Value = 25.69191589127711
$ws.Range("J5").Value = 9.826091555608782
$ws.Range("L5").Value = 13.07046429340689
$ws.Range("N5").Value = 18.45456975779411

# Row 6
$ws.Range("B6").Value = 17.35839585752686
$ws.Range("D6").Value = 7.989372413120847
$ws.Range("E6").Value = 12.76410539236529
$ws.Range("F6").Value = 37.13610054795073
$ws.Range("G6").Value = 44.39864182105831
$ws.Range("H6").Value = 17.77491111181981
$ws.Range("I6").Value = 25.69917426437953
$ws.Range("J6").Value = 9.824092497079532
$ws.Range("L6").Value = 13.06235823616969
$ws.Range("N6").Value = 18.4565414922259

# Row 7
$ws.Range("B7").Value = 17.45205657295864
$ws.Range("D7").Value = 8.014189259170221
$ws.Range("E7").Value = 12.80783693080043
$ws.Range("F7").Value = 37.15553726732238
$ws.Range("G7").Value = 44.53871650749713
$ws.Range("H7").Value = 17.77509015972172
$ws.Range("I7").Value = 25.64926249182583
$ws.Range("J7").Value = 9.837940900075591
$ws.Range("L7").Value = 13.11873479245529
$ws.Range("N7").Value = 18.44304009295782

# Row 8
$ws.Range("B8").Value = 17.86545789195446
$ws.Range("D8").Value = 8.121103502950264
$ws.Range("E8").Value = 12.99695425617334
$ws.Range("F8").Value = 37.26722808416613
$ws.Range("G8").Value = 45.19005327332934
$ws.Range("H8").Value = 17.78636236314325
$ws.Range("I8").Value = 25.44061582221426
$ws.Range("J8").Value = 9.898355686706156
$ws.Range("L8").Value = 13.36987961340567
$ws.Range("N8").Value = 18.38802245331356

# Row 9
$ws.Range("B9").Value = 18.66977226703289
$ws.Range("D9").Value = 8.320970593594312
$ws.Range("E9").Value = 13.3527607938316
$ws.Range("F9").Value = 37.5787341633706
$ws.Range("G9").Value = 46.57646908309327
$ws.Range("H9").Value = 17.84615211921223
$ws.Range("I9").Value = 25.07307953185558
$ws.Range("J9").Value = 10.01413155991764
$ws.Range("L9").Value = 13.86712703938943
$ws.Range("N9").Value = 18.29646402108877

# Row 10
$ws.Range("B10").Value = 19.25107117294269
$ws.Range("D10").Value = 8.461453600970524
$ws.Range("E10").Value = 13.60398517613219
$ws.Range("F10").Value = 37.86304998927843
$ws.Range("G10").Value = 47.65174284703573
$ws.Range("H10").Value = 17.91298957833267
$ws.Range("I10").Value = 24.82834299281459
$ws.Range("J10").Value = 10.0973731735391
$ws.Range("L10").Value = 14.23203618218331
$ws.Range("N10").Value = 18.23909643483699

# Row 11
$ws.Range("B11").Value = 19.51228179054754
$ws.Range("D11").Value = 8.523902451951789
$ws.Range("E11").Value = 13.71585601518085
$ws.Range("F11").Value = 38.00416014866714
$ws.Range("G11").Value = 48.15102060339342
$ws.Range("H11").Value = 17.94833668697557
$ws.Range("I11").Value = 24.72248291549465
$ws.Range("J11").Value = 10.13480958751722
$ws.Range("L11").Value = 14.39726905666481
$ws.Range("N11").Value = 18.21513674568

# Row 12
$ws.Range("B12").Value = 19.61064187849013
$ws.Range("D12").Value = 8.547334012436064
$ws.Range("E12").Value = 13.75785491777325
$ws.Range("F12").Value = 38.05925931368659
$ws.Range("G12").Value = 48.34135362738668
$ws.Range("H12").Value = 17.96242791571753
$ws.Range("I12").Value = 24.68318254644672
$ws.Range("J12").Value = 10.14892057998708
$ws.Range("L12").Value = 14.45967328687346
$ws.Range("N12").Value = 18.20637026077617

# Row 13
$ws.Range("B13").Value = 19.58948434084045
$ws.Range("D13").Value = 8.542297341168844
$ws.Range("E13").Value = 13.74882616742151
$ws.Range("F13").Value = 38.04731922427659
$ws.Range("G13").Value = 48.30030865281651
$ws.Range("H13").Value = 17.95936180054249
$ws.Range("I13").Value = 24.69161160562702
$ws.Range("J13").Value = 10.14588448728127
$ws.Range("L13").Value = 14.446241662717
$ws.Range("N13").Value = 18.20824465953273

# Row 14
$ws.Range("B14").Value = 19.52038555501044
$ws.Range("D14").Value = 8.525834557418706
$ws.Range("E14").Value = 13.71931868526026
$ws.Range("F14").Value = 38.00866002392473
$ws.Range("G14").Value = 48.1666551003899
$ws.Range("H14").Value = 17.94948186083386
$ws.Range("I14").Value = 24.71923388819374
$ws.Range("J14").Value = 10.13597183710808
$ws.Range("L14").Value = 14.40240672672562
$ws.Range("N14").Value = 18.21440938283769

# Row 15
$ws.Range("B15").Value = 19.47798569794349
$ws.Range("D15").Value = 8.515722236605971
$ws.Range("E15").Value = 13.70119657390326
$ws.Range("F15").Value = 37.98519590138444
$ws.Range("G15").Value = 48.08494786880098
$ws.Range("H15").Value = 17.94352190533314
$ws.Range("I15").Value = 24.73625575859308
$ws.Range("J15").Value = 10.12989143793819
$ws.Range("L15").Value = 14.37553334100607
$ws.Range("N15").Value = 18.21822534995734

# Row 16
$ws.Range("B16").Value = 19.23392766452794
$ws.Range("D16").Value = 8.4573425001216
$ws.Range("E16").Value = 13.59662412239759
$ws.Range("F16").Value = 37.85406260099905
$ws.Range("G16").Value = 47.61930153949388
$ws.Range("H16").Value = 17.91077863898603
$ws.Range("I16").Value = 24.83537127652481
$ws.Range("J16").Value = 10.09491765657462
$ws.Range("L16").Value = 14.22121760341699
$ws.Range("N16").Value = 18.24070519366562

# Row 17
$ws.Range("B17").Value = 19.08331372446585
$ws.Range("D17").Value = 8.421150155809968
$ws.Range("E17").Value = 13.53184198650111
$ws.Range("F17").Value = 37.77661232808812
$ws.Range("G17").Value = 47.33610264194153
$ws.Range("H17").Value = 17.89195456623945
$ws.Range("I17").Value = 24.89757664781699
$ws.Range("J17").Value = 10.07334963607229
$ws.Range("L17").Value = 14.12631344257897
$ws.Range("N17").Value = 18.25504265306029

# Row 18
$ws.Range("B18").Value = 18.99638614681206
$ws.Range("D18").Value = 8.400196481604361
$ws.Range("E18").Value = 13.49435494784939
$ws.Range("F18").Value = 37.73317512034298
$ws.Range("G18").Value = 47.17417904743054
$ws.Range("H18").Value = 17.88159300288061
$ws.Range("I18").Value = 24.93387057406272
$ws.Range("J18").Value = 10.06090395624559
$ws.Range("L18").Value = 14.07165818810553
$ws.Range("N18").Value = 18.26349038189626

# Row 19
$ws.Range("B19").Value = 18.96690555655881
$ws.Range("D19").Value = 8.393078645334416
$ws.Range("E19").Value = 13.48162418824461
$ws.Range("F19").Value = 37.71865953046228
$ws.Range("G19").Value = 47.1195259898514
$ws.Range("H19").Value = 17.87816483954363
$ws.Range("I19").Value = 24.94624754371855
$ws.Range("J19").Value = 10.05668324190212
$ws.Range("L19").Value = 14.05314275424419
$ws.Range("N19").Value = 18.26638521836262

# Row 20
$ws.Range("B20").Value = 19.09937839497907
$ws.Range("D20").Value = 8.425017100770491
$ws.Range("E20").Value = 13.53876167934442
$ws.Range("F20").Value = 37.78474234744905
$ws.Range("G20").Value = 47.36615115338552
$ws.Range("H20").Value = 17.89391026945747
$ws.Range("I20").Value = 24.89090147750538
$ws.Range("J20").Value = 10.07564978707057
$ws.Range("L20").Value = 14.13642369211948
$ws.Range("N20").Value = 18.25349558818406

# Row 21
$ws.Range("B21").Value = 19.54069729870122
$ws.Range("D21").Value = 8.530676004583869
$ws.Range("E21").Value = 13.72799576589842
$ws.Range("F21").Value = 38.01997026195698
$ws.Range("G21").Value = 48.20587950146474
$ws.Range("H21").Value = 17.95236471875414
$ws.Range("I21").Value = 24.71109920934649
$ws.Range("J21").Value = 10.13888522420072
$ws.Range("L21").Value = 14.41528704601741
$ws.Range("N21").Value = 18.2125903412394

# Row 22
$ws.Range("B22").Value = 19.82585273440147
$ws.Range("D22").Value = 8.59846544594078
$ws.Range("E22").Value = 13.8495409963963
$ws.Range("F22").Value = 38.18338742115969
$ws.Range("G22").Value = 48.76199023275429
$ws.Range("H22").Value = 17.99468060784865
$ws.Range("I22").Value = 24.5981725712501
$ws.Range("J22").Value = 10.17983064178156
$ws.Range("L22").Value = 14.59654911712611
$ws.Range("N22").Value = 18.187642730242

# Row 23
$ws.Range("B23").Value = 19.67398808570496
$ws.Range("D23").Value = 8.562402887479378
$ws.Range("E23").Value = 13.78487048647307
$ws.Range("F23").Value = 38.09529318551416
$ws.Range("G23").Value = 48.46457803023125
$ws.Range("H23").Value = 17.97172131729909
$ws.Range("I23").Value = 24.65802422199751
$ws.Range("J23").Value = 10.15801344768603
$ws.Range("L23").Value = 14.49991472630289
$ws.Range("N23").Value = 18.20079454621153

# Row 24
$ws.Range("B24").Value = 19.09211660209206
$ws.Range("D24").Value = 8.423269309852808
$ws.Range("E24").Value = 13.53563404155485
$ws.Range("F24").Value = 37.78106336845888
$ws.Range("G24").Value = 47.35256343103647
$ws.Range("H24").Value = 17.89302466087817
$ws.Range("I24").Value = 24.89391766731789
$ws.Range("J24").Value = 10.07461003162863
$ws.Range("L24").Value = 14.13185313556205
$ws.Range("N24").Value = 18.25419437783622

# Row 25
$ws.Range("B25").Value = 18.45347235008077
$ws.Range("D25").Value = 8.267990314673861
$ws.Range("E25").Value = 13.25822301236109
$ws.Range("F25").Value = 37.48464384374919
$ws.Range("G25").Value = 46.19071824972472
$ws.Range("H25").Value = 17.825945696638
$ws.Range("I25").Value = 25.16806087693466
$ws.Range("J25").Value = 9.983117803613165
$ws.Range("L25").Value = 13.73244541059914
$ws.Range("N25").Value = 18.29646402108877
